$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Rows 19-21 correspond to "ConQat (2nd)" RQ2/RQ3/RQ4, columns 3/4 are
# Ctags / BrlCad execution times that are currently empty.
$updates = @(
    @{ Row = 19; Col = 3; Text = "13 h 14 m" },
    @{ Row = 19; Col = 4; Text = "6 h 23 m" },
    @{ Row = 20; Col = 3; Text = "13 h 3 m" },
    @{ Row = 20; Col = 4; Text = "5 h 42 m" },
    @{ Row = 21; Col = 3; Text = "6 h 53 m" },
    @{ Row = 21; Col = 4; Text = "36 m 54 s" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
